$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header column in H1, copying the formatting from the
# neighboring header cell (G1) so it matches the other header cells' style.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Populate the data rows of the new "Save" column
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1
